$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with latest scraped values.
# D-column values are forced to remain text (matching the original inline-string
# cell type) even though some look numeric, then the style is reset to Normal so
# no extraneous number-format is left behind on the cell.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "44.544.15"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.61%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.427.02"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.61%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.41"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.57%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "101.67"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.09%  "
$ws.Range("E7").Value = "  +1.91%  "
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.515"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +5.61%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.27"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.91%  "
$ws.Range("E11").Value = "  +1.91%  "
$ws.Range("E12").Value = "  +1.24%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.88"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.47%  "
$ws.Range("E14").Value = "  +3.12%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.803.29"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.70%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.383.67"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.98%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.835"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +5.14%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "44.416.82"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.40%  "
$ws.Range("E19").Value = "  +4.23%  "
$ws.Range("E20").Value = "  +2.42%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0910"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.91"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.44%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "241.08"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.58%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.29"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.35%  "
$ws.Range("E25").Value = "  +1.81%  "
$ws.Range("E26").Value = "  -0.06%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.21"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.46%  "
$ws.Range("E28").Value = "  -4.27%  "
$ws.Range("E29").Value = "  +3.64%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "33.32"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.82%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "48.44"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.01%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.123"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +18.17%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.56"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +13.19%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.19"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.45%  "
$ws.Range("E35").Value = "  +0.28%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0769"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +7.14%  "
$ws.Range("E37").Value = "  +3.17%  "
$ws.Range("E38").Value = "  +4.02%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.90"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.47%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "124.57"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.40%  "
$ws.Range("E41").Value = "  +1.12%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.18"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.66%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.51"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.61%  "
$ws.Range("E44").Value = "  +3.71%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.947.25"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.54%  "
$ws.Range("E46").Value = "  +2.11%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.95"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +8.87%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.54"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.02%  "
$ws.Range("E49").Value = "  +10.49%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "53.57"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.01%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "73.62"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.09%  "
